$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "L474792"
$ws.Range("C5").Value = "NICK SHELL SERVICE"
$ws.Range("E5").Value = 560
$ws.Range("H5").Value = 45095.04187195602
$ws.Range("J5").Value = "06/13/23 11:00"
$ws.Range("K5").Value = "06/13/23 11:00"
$ws.Range("L5").Value = 60
$ws.Range("M5").Value = "`$560 as of 6/13/2023 9:00:01 AM"
$ws.Range("N5").Value = 700
$ws.Range("A6").Value = "LK644532"
$ws.Range("C6").Value = "SCL ENTERPRISES LAUNDRY"
$ws.Range("E6").Value = 880
$ws.Range("H6").Value = 45106.04187195602
$ws.Range("J6").Value = "06/12/23 19:34"
$ws.Range("K6").Value = "06/12/23 19:34"
$ws.Range("M6").Value = "`$880 as of 6/12/2023 5:34:00 PM"
$ws.Range("N6").Value = 920
$ws.Range("A7").Value = "L647934"
$ws.Range("C7").Value = "SB #6"
$ws.Range("E7").Value = 1940
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J7").Value = "04/06/23 22:10"
$ws.Range("K7").Value = "04/06/23 22:05"
$ws.Range("L7").Value = 20
$ws.Range("M7").Value = "`$1,940 as of 4/6/2023 8:05:45 PM"
$ws.Range("N7").Value = 1960
$ws.Range("A8").Value = "L688961"
$ws.Range("C8").Value = "MONA MART"
$ws.Range("E8").Value = 2640
$ws.Range("H8").Value = 45167.04187195602
$ws.Range("I8").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J8").Value = "06/13/23 14:33"
$ws.Range("K8").Value = "06/09/23 16:00"
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = "`$2,640 as of 6/9/2023 2:00:40 PM"
$ws.Range("N8").Value = 2640
$ws.Range("A9").Value = "LK561655"
$ws.Range("C9").Value = "CRENSHAW CRAVOR #2"
$ws.Range("E9").Value = 2780
$ws.Range("I9").Value = "ATM Inactive greater than 48 minutes"
$ws.Range("J9").Value = "01/23/20 08:24"
$ws.Range("K9").Value = "01/23/20 08:24"
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = "`$2,780 as of 1/23/2020 6:24:32 AM"
$ws.Range("N9").Value = 2800
$ws.Range("A10").Value = "L678988"
$ws.Range("C10").Value = "PAYELESS MARKET"
$ws.Range("E10").Value = 3000
$ws.Range("H10").Value = 45105.04187195602
$ws.Range("J10").Value = "06/12/23 15:37"
$ws.Range("K10").Value = "06/12/23 15:37"
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = "`$3,000 as of 6/12/2023 1:37:03 PM"
$ws.Range("N10").Value = 3000
$ws.Range("A11").Value = "L475182"
$ws.Range("C11").Value = "LA ESQUINA DE ORO"
$ws.Range("E11").Value = 3800
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = "ATM Inactive greater than 48 minutes"
$ws.Range("J11").Value = "09/16/20 16:57"
$ws.Range("K11").Value = "09/15/20 23:38"
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = "`$3,800 as of 9/16/2020 1:28:00 PM"
$ws.Range("N11").Value = 3800
$ws.Range("A12").Value = "L682801"
$ws.Range("C12").Value = "SB#5"
$ws.Range("E12").Value = 3880
$ws.Range("H12").Value = 45106.04187195602
$ws.Range("J12").Value = "06/13/23 16:25"
$ws.Range("K12").Value = "06/13/23 16:25"
$ws.Range("M12").Value = "`$4,020 as of 6/12/2023 6:05:36 PM"
$ws.Range("N12").Value = 3920
$ws.Range("A13").Value = "LK236828"
$ws.Range("C13").Value = "WORLDWIDE AUTOMOTIVE"
$ws.Range("E13").Value = 5320
$ws.Range("H13").Value = 45113.04187195602
$ws.Range("I13").ClearContents()
$ws.Range("J13").Value = "06/12/23 20:28"
$ws.Range("K13").Value = "06/12/23 20:28"
$ws.Range("L13").Value = 80
$ws.Range("M13").Value = "`$5,320 as of 6/12/2023 6:28:00 PM"
$ws.Range("N13").Value = 5360
$ws.Range("A14").Value = "L662336"
$ws.Range("C14").Value = "SB#4 MONA MARKET"
$ws.Range("E14").Value = 6080
$ws.Range("H14").Value = 45117.04187195602
$ws.Range("I14").ClearContents()
$ws.Range("J14").Value = "06/13/23 16:48"
$ws.Range("K14").Value = "06/13/23 15:07"
$ws.Range("M14").Value = "`$6,220 as of 6/13/2023 8:24:02 AM"
$ws.Range("N14").Value = 6080
$ws.Range("A15").Value = "L474817"
$ws.Range("C15").Value = "SAFETY MARKET"
$ws.Range("E15").Value = 6320
$ws.Range("H15").Value = 45099.04187195602
$ws.Range("J15").Value = "06/13/23 16:26"
$ws.Range("K15").Value = "06/13/23 15:26"
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = "`$6,440 as of 6/13/2023 7:36:23 AM"
$ws.Range("N15").Value = 6320
$ws.Range("A16").Value = "L488595"
$ws.Range("C16").Value = "N S MART"
$ws.Range("E16").Value = 6480
$ws.Range("H16").Value = 45128.04187195602
$ws.Range("I16").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J16").Value = "06/12/23 15:05"
$ws.Range("K16").Value = "06/11/23 21:53"
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = "`$6,480 as of 6/11/2023 7:53:56 PM"
$ws.Range("N16").Value = 6480
$ws.Range("A17").Value = "L476340"
$ws.Range("C17").Value = "DONUT & SANDWICH"
$ws.Range("E17").Value = 6620
$ws.Range("H17").Value = 45117.04187195602
$ws.Range("I17").ClearContents()
$ws.Range("J17").Value = "06/13/23 14:15"
$ws.Range("K17").Value = "06/13/23 14:15"
$ws.Range("L17").Value = 40
$ws.Range("M17").Value = "`$6,660 as of 6/13/2023 5:11:45 AM"
$ws.Range("N17").Value = 6660
$ws.Range("A18").Value = "L474746"
$ws.Range("C18").Value = "ZACATES MARKET"
$ws.Range("E18").Value = 6780
$ws.Range("H18").Value = 45115.04187195602
$ws.Range("J18").Value = "06/12/23 19:10"
$ws.Range("K18").Value = "06/12/23 19:10"
$ws.Range("M18").Value = "`$6,780 as of 6/12/2023 5:10:31 PM"
$ws.Range("N18").Value = 6840
$ws.Range("A19").Value = "L697590"
$ws.Range("C19").Value = "S B MARKET ST"
$ws.Range("E19").Value = 6900
$ws.Range("H19").Value = 45103.04187195602
$ws.Range("J19").Value = "06/13/23 16:35"
$ws.Range("K19").Value = "06/13/23 16:35"
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = "`$7,060 as of 6/13/2023 4:03:51 AM"
$ws.Range("N19").Value = 6920
$ws.Range("A20").Value = "LK864765"
$ws.Range("C20").Value = "SKY LIQUOR"
$ws.Range("E20").Value = 7120
$ws.Range("H20").Value = 45103.04187195602
$ws.Range("J20").Value = "06/13/23 16:41"
$ws.Range("K20").Value = "06/13/23 16:41"
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = "`$7,340 as of 6/13/2023 11:29:34 AM"
$ws.Range("N20").Value = 7120
$ws.Range("A21").Value = "L474761"
$ws.Range("C21").Value = "BABS MARKET"
$ws.Range("E21").Value = 7180
$ws.Range("H21").Value = 45156.04187195602
$ws.Range("J21").Value = "06/12/23 20:20"
$ws.Range("K21").Value = "06/12/23 20:20"
$ws.Range("L21").Value = 40
$ws.Range("M21").Value = "`$7,180 as of 6/12/2023 6:20:05 PM"
$ws.Range("N21").Value = 7220
$ws.Range("A22").Value = "L688966"
$ws.Range("C22").Value = "LACON MINI MART"
$ws.Range("E22").Value = 7300
$ws.Range("H22").Value = 45165.04187195602
$ws.Range("I22").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J22").Value = "06/11/23 15:48"
$ws.Range("K22").Value = "06/11/23 15:48"
$ws.Range("M22").Value = "`$7,300 as of 6/11/2023 1:48:15 PM"
$ws.Range("N22").Value = 7400
$ws.Range("A23").Value = "L697589"
$ws.Range("C23").Value = "S B DISCOUNT MART"
$ws.Range("E23").Value = 7740
$ws.Range("H23").Value = 45096.04187195602
$ws.Range("J23").Value = "06/13/23 10:33"
$ws.Range("K23").Value = "06/13/23 10:33"
$ws.Range("L23").Value = 20
$ws.Range("M23").Value = "`$7,740 as of 6/13/2023 8:33:20 AM"
$ws.Range("N23").Value = 7780
$ws.Range("A24").Value = "LK923383"
$ws.Range("C24").Value = "SAMYS PHONE CARDS"
$ws.Range("E24").Value = 10220
$ws.Range("H24").Value = 45104.04187195602
$ws.Range("I24").ClearContents()
$ws.Range("J24").Value = "06/12/23 22:27"
$ws.Range("K24").Value = "06/12/23 22:27"
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = "`$10,220 as of 6/12/2023 8:27:39 PM"
$ws.Range("N24").Value = 10320
$ws.Range("A25").Value = "L475090"
$ws.Range("C25").Value = "S.B. 2"
$ws.Range("E25").Value = 13420
$ws.Range("H25").Value = 45107.04187195602
$ws.Range("J25").Value = "06/13/23 16:30"
$ws.Range("K25").Value = "06/13/23 14:25"
$ws.Range("M25").Value = "`$13,620 as of 6/13/2023 11:54:34 AM"
$ws.Range("N25").Value = 13320
$ws.Range("E26").Value = 23220
$ws.Range("H26").Value = 45103.04187195602
$ws.Range("J26").Value = "06/13/23 16:38"
$ws.Range("K26").Value = "06/13/23 16:38"
$ws.Range("L26").Value = 60
$ws.Range("M26").Value = "`$23,440 as of 6/13/2023 11:03:04 AM"
$ws.Range("N26").Value = 23420
$ws.Range("E27").Value = 140180
